$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns L and M entirely (Likely Closed / Escrow Officer)
$ws.Range("L1:M126").Delete()

# Clear the GF# data values (column K) for rows 2-126, keep header in K1
$ws.Range("K2:K126").Clear()

$ws.Range("K3").Select()
